$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item(1).Name = "GNG_TO-16512555442331362"
$wb.Worksheets.Item(2).Name = "NB_TO-16512555465662801"
$wb.Worksheets.Item(3).Name = "RS_TO-16512555465722828"
$wb.Worksheets.Item(4).Name = "TOL_TO-16512555466312845"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16512555467092834"

# Sheet1 (GNG_TO) B2:B5
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16512555442001367.csv"
$ws1.Range("B3").Value = "GNG_stims-1651255544215136.csv"
$ws1.Range("B4").Value = "go_stims-16512555442171378.csv"
$ws1.Range("B5").Value = "GNG_stims-16512555442321372.csv"

# Sheet2 (NB_TO) B2:B10
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16512555450671253.csv"
$ws2.Range("B3").Value = "OB-16512555450381255.csv"
$ws2.Range("B4").Value = "ZB-match_2-16512555444861372.csv"
$ws2.Range("B5").Value = "ZB-match_1-16512555443201363.csv"
$ws2.Range("B6").Value = "TB-1651255546548289.csv"
$ws2.Range("B7").Value = "TB-1651255546296281.csv"
$ws2.Range("B8").Value = "OB-1651255544752129.csv"
$ws2.Range("B9").Value = "ZB-match_8-16512555443991375.csv"
$ws2.Range("B10").Value = "TB-16512555455212822.csv"

# Sheet4 (TOL_TO) B2:B7
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-1651255546598283.csv"
$ws4.Range("B3").Value = "ZM_stims-1651255546575288.csv"
$ws4.Range("B4").Value = "MM_stims-16512555466142814.csv"
$ws4.Range("B5").Value = "ZM_stims-16512555465993447.csv"
$ws4.Range("B6").Value = "MM_stims-16512555466302814.csv"
$ws4.Range("B7").Value = "ZM_stims-16512555466152825.csv"

# Sheet5 (vSAT_TO) B2:B5
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16512555466362832.csv"
$ws5.Range("B3").Value = "vSAT_stims-16512555466942828.csv"
$ws5.Range("B4").Value = "vSAT_stims-16512555466772816.csv"
$ws5.Range("B5").Value = "SAT_stims-16512555466612833.csv"
